$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 3 (3rd data row): John Gerega / 2/13/26 / two-sentence task
$tbl.Cell(3, 1).Range.Text = "John Gerega"
$tbl.Cell(3, 2).Range.Text = "2/13/26"

$taskXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Fixed buffer issue on the program. </w:t></w:r><w:r><w:t>Testing to ensure scanner still worked properly</w:t></w:r></w:p>'
[void]$tbl.Cell(3, 3).Range.InsertXML($taskXml)

# Row 4 (4th data row): John Gerega / 2/16/26 / Comment through new changes
$tbl.Cell(4, 1).Range.Text = "John Gerega"
$tbl.Cell(4, 2).Range.Text = "2/16/26"
$tbl.Cell(4, 3).Range.Text = "Comment through new changes"

# Row 5 (5th data row): John Gerega / 2/18/26 / Testing with different input files
$tbl.Cell(5, 1).Range.Text = "John Gerega"
$tbl.Cell(5, 2).Range.Text = "2/18/26"
$tbl.Cell(5, 3).Range.Text = "Testing with different input files"
